$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row and title-case municipality/state names (standalone
# "de/del/la/las/el/los/y" connector words capitalized), plus one float
# precision correction in D449.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B10").Value = "Playas De Rosarito"
$ws.Range("B26").Value = "Chiapa De Corzo"
$ws.Range("B29").Value = "Comitán De Domínguez"
$ws.Range("B45").Value = "Marqués De Comillas"
$ws.Range("B46").Value = "Mazapa De Madero"
$ws.Range("B49").Value = "Ocozocoautla De Espinosa"
$ws.Range("B53").Value = "San Cristóbal De Las Casas"
$ws.Range("B79").Value = "Hidalgo Del Parral"
$ws.Range("B85").Value = "San Francisco Del Oro"
$ws.Range("B107").Value = "Villa De Álvarez"
$ws.Range("A109").Value = "Ciudad De México"
$ws.Range("B125").Value = "Coneto De Comonfort"
$ws.Range("B137").Value = "Nombre De Dios"
$ws.Range("B140").Value = "Pánuco De Coronado"
$ws.Range("B144").Value = "San Juan Del Río"
$ws.Range("A152").Value = "Estado De México"
$ws.Range("B152").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B154").Value = "Almoloya De Juárez"
$ws.Range("B157").Value = "Atizapán De Zaragoza"
$ws.Range("B160").Value = "Chapa De Mota"
$ws.Range("B165").Value = "Ecatepec De Morelos"
$ws.Range("B169").Value = "Ixtapan De La Sal"
$ws.Range("B178").Value = "Naucalpan De Juárez"
$ws.Range("B181").Value = "San Felipe Del Progreso"
$ws.Range("B182").Value = "Soyaniquilpan De Juárez"
$ws.Range("B194").Value = "Tlalnepantla De Baz"
$ws.Range("B198").Value = "Valle De Bravo"
$ws.Range("B199").Value = "Villa De Allende"
$ws.Range("B207").Value = "San Miguel De Allende"
$ws.Range("B208").Value = "Apaseo El Alto"
$ws.Range("B209").Value = "Apaseo El Grande"
$ws.Range("B216").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B227").Value = "San Diego De La Unión"
$ws.Range("B229").Value = "San Francisco Del Rincón"
$ws.Range("B231").Value = "San Luis De La Paz"
$ws.Range("B232").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B234").Value = "Silao De La Victoria"
$ws.Range("B238").Value = "Valle De Santiago"
$ws.Range("B244").Value = "Acapulco De Juárez"
$ws.Range("B247").Value = "Ajuchitlán Del Progreso"
$ws.Range("B248").Value = "Alcozauca De Guerrero"
$ws.Range("B253").Value = "Atoyac De Álvarez"
$ws.Range("B254").Value = "Ayutla De Los Libres"
$ws.Range("B257").Value = "Buenavista De Cuéllar"
$ws.Range("B258").Value = "Chilapa De Álvarez"
$ws.Range("B259").Value = "Chilpancingo De Los Bravo"
$ws.Range("B260").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B264").Value = "Coyuca De Benítez"
$ws.Range("B265").Value = "Coyuca De Catalán"
$ws.Range("B268").Value = "Cutzamala De Pinzón"
$ws.Range("B273").Value = "Huitzuco De Los Figueroa"
$ws.Range("B274").Value = "Iguala De La Independencia"
$ws.Range("B275").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B276").Value = "Zihuatanejo De Azueta"
$ws.Range("B278").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B280").Value = "Mártir De Cuilapan"
$ws.Range("B291").Value = "Taxco De Alarcón"
$ws.Range("B293").Value = "Técpan De Galeana"
$ws.Range("B295").Value = "Tepecoacuilco De Trujano"
$ws.Range("B298").Value = "Tlapa De Comonfort"
$ws.Range("B310").Value = "Atotonilco De Tula"
$ws.Range("B311").Value = "Atotonilco El Grande"
$ws.Range("B315").Value = "Cuautepec De Hinojosa"
$ws.Range("B320").Value = "Huejutla De Reyes"
$ws.Range("B323").Value = "Jacala De Ledezma"
$ws.Range("B328").Value = "Mixquiahuala De Juárez"
$ws.Range("B329").Value = "Pachuca De Soto"
$ws.Range("B335").Value = "Santiago De Anaya"
$ws.Range("B340").Value = "Tepehuacán De Guerrero"
$ws.Range("B341").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B342").Value = "Tezontepec De Aldama"
$ws.Range("B348").Value = "Tulancingo De Bravo"
$ws.Range("B350").Value = "Zacualtipán De Ángeles"
$ws.Range("B354").Value = "Ahualulco De Mercado"
$ws.Range("B357").Value = "Atotonilco El Alto"
$ws.Range("B368").Value = "Encarnación De Díaz"
$ws.Range("B371").Value = "Ixtlahuacán Del Río"
$ws.Range("B375").Value = "Lagos De Moreno"
$ws.Range("B378").Value = "Ojuelos De Jalisco"
$ws.Range("B382").Value = "San Cristóbal De La Barranca"
$ws.Range("B383").Value = "San Juan De Los Lagos"
$ws.Range("B386").Value = "San Miguel El Alto"
$ws.Range("B387").Value = "San Sebastián Del Oeste"
$ws.Range("B389").Value = "Talpa De Allende"
$ws.Range("B390").Value = "Tamazula De Gordiano"
$ws.Range("B393").Value = "Tepatitlán De Morelos"
$ws.Range("B394").Value = "Tizapán El Alto"
$ws.Range("B395").Value = "Tlajomulco De Zúñiga"
$ws.Range("B402").Value = "Unión De San Antonio"
$ws.Range("B403").Value = "Valle De Guadalupe"
$ws.Range("B404").Value = "Valle De Juárez"
$ws.Range("B406").Value = "Yahualica De González Gallo"
$ws.Range("B409").Value = "Zapotlán El Grande"
$ws.Range("D449").Value = 0.009557945041816007
$ws.Range("B470").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B492").Value = "Coatlán Del Río"
$ws.Range("B498").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B500").Value = "Puente De Ixtla"
$ws.Range("B505").Value = "Tetela Del Volcán"
$ws.Range("B506").Value = "Tlaltizapán De Zapata"
$ws.Range("B512").Value = "Zacualpan De Amilpas"
$ws.Range("B517").Value = "Ixtlán Del Río"
$ws.Range("B542").Value = "San Nicolás De Los Garza"
$ws.Range("B547").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B549").Value = "Coicoyán De Las Flores"
$ws.Range("B551").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B552").Value = "Guelatao De Juárez"
$ws.Range("B553").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B554").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B555").Value = "Ixtlán De Juárez"
$ws.Range("B556").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B559").Value = "Mariscala De Juárez"
$ws.Range("B561").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B562").Value = "Nejapa De Madero"
$ws.Range("B563").Value = "Oaxaca De Juárez"
$ws.Range("B564").Value = "Ocotlán De Morelos"
$ws.Range("B565").Value = "Putla Villa De Guerrero"
$ws.Range("B604").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B630").Value = "Santo Domingo De Morelos"
$ws.Range("B635").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B636").Value = "Tlacolula De Matamoros"
$ws.Range("B638").Value = "Villa De Tututepec"
$ws.Range("B639").Value = "Villa De Zaachila"
$ws.Range("B640").Value = "Villa Sola De Vega"
$ws.Range("B641").Value = "Zimatlán De Álvarez"
$ws.Range("B659").Value = "Cuapiaxtla De Madero"
$ws.Range("B660").Value = "Cuayuca De Andrade"
$ws.Range("B669").Value = "Huehuetlán El Grande"
$ws.Range("B672").Value = "Ixcamilpa De Guerrero"
$ws.Range("B673").Value = "Izúcar De Matamoros"
$ws.Range("B679").Value = "Los Reyes De Juárez"
$ws.Range("B692").Value = "San Salvador El Verde"
$ws.Range("B700").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B703").Value = "Tepexi De Rodríguez"
$ws.Range("B708").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B716").Value = "Xayacatlán De Bravo"
$ws.Range("B727").Value = "Amealco De Bonfil"
$ws.Range("B729").Value = "Cadereyta De Montes"
$ws.Range("B733").Value = "Jalpan De Serra"
$ws.Range("B734").Value = "Landa De Matamoros"
$ws.Range("B737").Value = "Pinal De Amoles"
$ws.Range("B739").Value = "San Juan Del Río"
$ws.Range("B749").Value = "Axtla De Terrazas"
$ws.Range("B752").Value = "Cerro De San Pedro"
$ws.Range("B754").Value = "Ciudad Del Maíz"
$ws.Range("B762").Value = "Mexquitic De Carmona"
$ws.Range("B770").Value = "Santa María Del Río"
$ws.Range("B772").Value = "Soledad De Graciano Sánchez"
$ws.Range("B781").Value = "Villa De Arista"
$ws.Range("B782").Value = "Villa De Guadalupe"
$ws.Range("B783").Value = "Villa De Ramos"
$ws.Range("B784").Value = "Villa De Reyes"
$ws.Range("B839").Value = "Soto La Marina"
$ws.Range("B852").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B853").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B868").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B871").Value = "Amatlán De Los Reyes"
$ws.Range("B876").Value = "Camarón De Tejeda"
$ws.Range("B877").Value = "Castillo De Teayo"
$ws.Range("B884").Value = "Cosamaloapan De Carpio"
$ws.Range("B890").Value = "Hueyapan De Ocampo"
$ws.Range("B891").Value = "Ignacio De La Llave"
$ws.Range("B895").Value = "Ixhuatlán De Madero"
$ws.Range("B896").Value = "Ixhuatlán Del Café"
$ws.Range("B897").Value = "Ixhuatlán Del Sureste"
$ws.Range("B903").Value = "Juchique De Ferrer"
$ws.Range("B906").Value = "Lerdo De Tejada"
$ws.Range("B908").Value = "Martínez De La Torre"
$ws.Range("B917").Value = "Ozuluama De Mascareñas"
$ws.Range("B920").Value = "Paso De Ovejas"
$ws.Range("B921").Value = "Paso Del Macho"
$ws.Range("B923").Value = "Poza Rica De Hidalgo"
$ws.Range("B929").Value = "Soledad De Doblado"
$ws.Range("B953").Value = "Vega De Alatorre"
$ws.Range("B969").Value = "Concepción Del Oro"
$ws.Range("B981").Value = "Moyahua De Estrada"
$ws.Range("B982").Value = "Nochistlán De Mejía"
$ws.Range("B983").Value = "Noria De Ángeles"
$ws.Range("B991").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B993").Value = "Villa De Cos"

# Drop the trailing metadata/footer rows (sample size, source, author, date)
# so the sheet ends at row 999 -- dimension recalculates automatically.
$ws.Range("1001:1005").Delete()

